# "multiple static page tests added"
# Marks the newly-added "Static Page test" / "Multiple page test from array"
# rows as done (J9, J10 go from "?" to "x"), marks a few other previously
# unanswered rows as done too (J15, J16, J19, J20, J21: "?" -> "x"), bumps
# the repeatable "complex xpath" count (J11: 8 -> 10), grows row 2's height
# so the wrapped header text fits, and leaves the selection on J8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that were answered with a placeholder "?" and are now marked done ("x")
$ws.Range("J9").Value  = "x"
$ws.Range("J10").Value = "x"
$ws.Range("J15").Value = "x"
$ws.Range("J16").Value = "x"
$ws.Range("J19").Value = "x"
$ws.Range("J20").Value = "x"
$ws.Range("J21").Value = "x"

# "complex xpath" repeatable count bumped from 8 to 10
$ws.Range("J11").Value = 10

# Row 2 (wrapped header row) grows taller to fit its content
$ws.Rows.Item(2).RowHeight = 36

# Leave the active selection on J8
$ws.Range("J8").Select() | Out-Null
